$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("Z2:Z4").Value = "2025-11-13T06:52:41.035320"
$ws.Range("Z5:Z12").Value = "2025-11-13T06:52:41.036322"
$ws.Range("Z13").Value = "2025-11-13T06:52:41.037509"
$ws.Range("Z14:Z15").Value = "2025-11-13T06:52:41.037851"
$ws.Range("Z16:Z19").Value = "2025-11-13T06:52:41.038239"
$ws.Range("Z20:Z25").Value = "2025-11-13T06:52:41.039038"
$ws.Range("Z26:Z34").Value = "2025-11-13T06:52:41.039938"
$ws.Range("Z35:Z44").Value = "2025-11-13T06:52:41.040943"
$ws.Range("Z45").Value = "2025-11-13T06:52:41.041940"
$ws.Range("Z46:Z47").Value = "2025-11-13T06:52:41.381670"
$ws.Range("Z48:Z51").Value = "2025-11-13T06:52:41.382218"
$ws.Range("Z52:Z60").Value = "2025-11-13T06:52:41.382742"
$ws.Range("Z61:Z62").Value = "2025-11-13T06:52:41.383738"
$ws.Range("Z63:Z67").Value = "2025-11-13T06:52:41.384005"
$ws.Range("Z68").Value = "2025-11-13T06:52:41.384514"
$ws.Range("Z69:Z71").Value = "2025-11-13T06:52:41.384560"
$ws.Range("Z72:Z73").Value = "2025-11-13T06:52:41.385102"
$ws.Range("Z74").Value = "2025-11-13T06:52:41.385643"
$ws.Range("Z75:Z80").Value = "2025-11-13T06:52:41.640873"
$ws.Range("Z81:Z89").Value = "2025-11-13T06:52:41.641863"
$ws.Range("Z90:Z99").Value = "2025-11-13T06:52:41.642860"
$ws.Range("Z100:Z102").Value = "2025-11-13T06:52:41.643860"
